$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value as literal TEXT (avoid Excel auto-converting
# percent-looking / numeric-looking / boolean-looking strings into a
# number/bool) and make sure no residual number-format/style (quote
# prefix) is left behind on the cell.
function Set-TextValue($addr, $val) {
    $r = $ws.Range($addr)
    $r.Value = "'" + $val
    $r.ClearFormats()
}

function Set-BoolValue($addr, $val) {
    $ws.Range($addr).Value = $val
}

# ---------------------------------------------------------------
# Row 2 - Dom Pérignon Vintage Champagne3 -> ...Champagne4
# ---------------------------------------------------------------
$ws.Range("C2").Value = "Dom Pérignon Vintage Champagne4"
$ws.Range("D2").Value = "Dom Pérignon Vintage Champagne4"
$ws.Range("AU2").Value = "2025-03-29T09:07:10.002Z"

# ---------------------------------------------------------------
# Row 3 - Macallan Rare Cask Single Malt2 -> Macallan Rare Cask Single Malt
# ---------------------------------------------------------------
$ws.Range("C3").Value = "Macallan Rare Cask Single Malt"
$ws.Range("D3").Value = "Macallan Rare Cask Single Malt"
$ws.Range("AC3").Value = "750 ML"
$ws.Range("AU3").Value = "2025-03-28T15:21:02.207Z"

# ---------------------------------------------------------------
# Row 16 - Casamigos Blanco Tequila3 -> Casamigos Blanco Tequila
# ---------------------------------------------------------------
$ws.Range("C16").Value = "Casamigos Blanco Tequila"
$ws.Range("D16").Value = "Casamigos Blanco Tequila"

Set-BoolValue "W16" $true
Set-BoolValue "X16" $false
Set-BoolValue "Y16" $false
Set-BoolValue "Z16" $true

Set-TextValue "AA16" "40%"

$ws.Range("AC16").Value = "750 ML"

Set-BoolValue "AS16" $true

# AT16 / AU16 removed entirely
$ws.Range("AT16").ClearContents()
$ws.Range("AU16").ClearContents()

# ---------------------------------------------------------------
# New row 26
# ---------------------------------------------------------------
Set-TextValue "A26" "608"
$ws.Range("B26").Value = "SKU-608-717"
$ws.Range("C26").Value = "test product"
$ws.Range("D26").Value = "test product"
$ws.Range("F26").Value = "beer"
$ws.Range("G26").Value = "whiskey"
$ws.Range("Q26").Value = 34
$ws.Range("S26").Value = 34
$ws.Range("T26").Value = "https://res.cloudinary.com/dc3hqcovg/image/upload/v1743226484/vzcjhlpqnvz4i0x2uw1j.svg"
Set-TextValue "W26" "false"
Set-TextValue "X26" "false"
Set-TextValue "Y26" "false"
Set-TextValue "Z26" "false"
$ws.Range("AB26").Value = "45ml"
$ws.Range("AC26").Value = "45ML"
$ws.Range("AG26").Value = 34
$ws.Range("AK26").Value = "taxable"
$ws.Range("AL26").Value = "no"
$ws.Range("AM26").Value = "no"
$ws.Range("AN26").Value = "instock"
$ws.Range("AP26").Value = "publish"
$ws.Range("AQ26").Value = "open"
Set-TextValue "AS26" "false"
$ws.Range("AT26").Value = "https://res.cloudinary.com/dc3hqcovg/image/upload/v1743226484/vzcjhlpqnvz4i0x2uw1j.svg"
$ws.Range("AU26").Value = "2025-03-29T05:34:46.251Z"
$ws.Range("AV26").Value = "beer"
Set-TextValue "AW26" "608"
$ws.Range("AX26").Value = "2025-03-29T05:34:46.251Z"

# ---------------------------------------------------------------
# New row 27
# ---------------------------------------------------------------
Set-TextValue "A27" "609"
$ws.Range("B27").Value = "SKU-609-342"
$ws.Range("C27").Value = "waefds"
$ws.Range("D27").Value = "waefds"
$ws.Range("F27").Value = "test category 34"
$ws.Range("G27").Value = "test"
$ws.Range("H27").Value = "test"
$ws.Range("N27").Value = "test"
$ws.Range("Q27").Value = 25
$ws.Range("S27").Value = 25
$ws.Range("T27").Value = "https://res.cloudinary.com/dc3hqcovg/image/upload/v1743239260/awsdfs6a1sfgwh2b3ik6.jpg"
Set-TextValue "W27" "false"
Set-TextValue "X27" "false"
Set-TextValue "Y27" "false"
Set-TextValue "Z27" "false"
$ws.Range("AB27").Value = "ml"
$ws.Range("AC27").Value = "ML"
$ws.Range("AG27").Value = 23
$ws.Range("AI27").Value = "test"
$ws.Range("AK27").Value = "taxable"
$ws.Range("AL27").Value = "no"
$ws.Range("AM27").Value = "no"
$ws.Range("AN27").Value = "instock"
$ws.Range("AP27").Value = "publish"
$ws.Range("AQ27").Value = "open"
Set-TextValue "AS27" "false"
$ws.Range("AT27").Value = "https://res.cloudinary.com/dc3hqcovg/image/upload/v1743239260/awsdfs6a1sfgwh2b3ik6.jpg"
$ws.Range("AU27").Value = "2025-03-29T09:07:41.265Z"
$ws.Range("AV27").Value = "test category 34"
Set-TextValue "AW27" "609"
$ws.Range("AX27").Value = "2025-03-29T09:07:41.265Z"
